# Keenan_Board_2.0 Parts_List.xlsx -- "Added switch to parts list"
#
# Inserts a new "Power Switch" line between the "Dual 4:1 Muxes" row and the
# "TI Stellaris EK-LM4F120XL" row, and adds a Cost/Unit total (SUM) under the
# existing "Total:" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Push row 5 ("TI Stellaris...") and row 6 ("Total:") down one row, opening
# up a blank row 5 for the new part. This carries formatting, formulas and
# row heights along with it.
$ws.Rows("5:5").Insert()

# --- Fill in the new row 5: Power Switch -----------------------------------
# (The link text is written first so the two new shared-string entries land
# in the same order as the authored workbook: url then part name.)
$ws.Range("E5").Value = "https://www.sparkfun.com/products/9609"
$ws.Range("A5").Value = "Power Switch"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0.75
$ws.Range("C5").NumberFormat = $ws.Range("C6").NumberFormat

# --- Re-establish column D's formula down through the shifted/new rows -----
# (Writing D5/D6 explicitly also nudges D3/D4 back onto the shared formula
# group that "Rows.Insert()" otherwise breaks apart.)
$ws.Range("D5").Formula = "=(C5*B5)"
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("D6").Formula = "=(C6*B6)"

# --- Cost/Unit total under the parts list -----------------------------------
$ws.Range("C7").Formula = "=SUM(C2:C6)"
$ws.Range("C7").NumberFormat = $ws.Range("C6").NumberFormat

# --- Hyperlinks: rebuild clean so the shifted row keeps pointing at the
#     TI Stellaris / Mouser link, and the new row gets the Sparkfun link. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "http://www.newark.com/stmicroelectronics/ld1117s33ctr/ic-ldo-volt-reg-3-3v-0-8a-sot/dp/89K0626?CMP=AFC-OP")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.sparkfun.com/products/9473")
$ws.Hyperlinks.Add($ws.Range("E4"), "http://www.newark.com/nxp/74hc4052d-653/ic-analog-mux-dmux-dual-4-x-1/dp/78R7402")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.sparkfun.com/products/9609")
$ws.Hyperlinks.Add($ws.Range("E6"), "http://www.mouser.com/ProductDetail/Texas-Instruments/EK-LM4F120XL/?qs=t9Lg9qrXjEy2enepSwqR9A==")

# Hyperlinks.Add() mints a brand-new "Hyperlink"-ish style for the cell
# instead of reusing the sheet's existing Hyperlink style (xfId 2, s="2").
# Re-apply the named style so E2:E6 stay on the original style index.
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E6").Style = "Hyperlink"

# --- Selection matches the saved view in the authored workbook -------------
$ws.Range("C8").Select()
